# Agenda de projeto - add new task "Trigges Controle de estoque" rows (12 and 13)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Row 12: "Criar" task for the new entry
$ws.Range("A12").Value = "Criar"
$ws.Range("B12").Value = "Trigges Controle de estoque"
$ws.Range("C12").Value = "T4 - Trigges"
$ws.Range("D12").Value = 9
$ws.Range("E12").Value = "Neimar"
$ws.Range("F12").Value = "?"
$ws.Range("G12").Value = "Notepad++"

# Row 13: "Teste" task for the new entry
$ws.Range("A13").Value = "Teste"
$ws.Range("B13").Value = "Trigges Controle de estoque"
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = "Aurélio"
$ws.Range("F13").Value = "?"
$ws.Range("G13").Value = "Postgree"

# Fix B13 style (remove border so it matches row 12's plain style)
$ws.Range("B13").Borders.LineStyle = -4142

# Update the selection to reflect where the user left off editing
$ws.Range("H13").Select()

$wb.Save()
